$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.035.94"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").Value = "2.300.30"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.81%  "

$ws.Range("E7").Value = "  -0.47%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.512"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.62%  "

$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("E12").Value = "  +0.63%  "

$ws.Range("E13").Value = "  +4.12%  "

$ws.Range("E14").Value = "  +1.91%  "

$ws.Range("D15").Value = "2.656.89"
$ws.Range("E15").Value = "  -0.27%  "

$ws.Range("D16").Value = "2.330.87"
$ws.Range("E16").Value = "  +1.18%  "

$ws.Range("E17").Value = "  -2.20%  "

$ws.Range("D18").Value = "42.957.97"
$ws.Range("E18").Value = "  -0.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.72%  "

$ws.Range("E20").Value = "  +0.31%  "

$ws.Range("E21").Value = "  +0.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.20%  "

$ws.Range("E25").Value = "  +0.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("E30").Value = "  -4.82%  "

$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("E32").Value = "  -0.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.36%  "

$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("E37").Value = "  -1.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.102"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("E39").Value = "  -0.72%  "

$ws.Range("E40").Value = "  +0.81%  "

$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("E42").Value = "  +3.09%  "

$ws.Range("E43").Value = "  -3.50%  "

$ws.Range("D44").Value = "1.982.20"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.44%  "

$ws.Range("E46").Value = "  +1.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.27%  "

$ws.Range("E49").Value = "  +3.79%  "

$ws.Range("D50").Value = "2.523.78"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("E51").Value = "  +0.56%  "

